$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "row/column index" label cells A6:D9 (text, not numbers) ---
$ws.Range("A6:D9").NumberFormat = "@"

$ws.Range("A6").Value = "00"
$ws.Range("B6").Value = "01"
$ws.Range("C6").Value = "02"
$ws.Range("D6").Value = "03"

$ws.Range("A7").Value = "10"
$ws.Range("B7").Value = "11"
$ws.Range("C7").Value = "12"
$ws.Range("D7").Value = "13"

$ws.Range("A8").Value = "20"
$ws.Range("B8").Value = "21"
$ws.Range("C8").Value = "22"
$ws.Range("D8").Value = "23"

$ws.Range("A9").Value = "30"
$ws.Range("B9").Value = "31"
$ws.Range("C9").Value = "32"
$ws.Range("D9").Value = "33"

# --- New "last added number" highlighting formulas, rows 31-34, cols F-I ---
# Row 31 references row 6 (indices) / row 1 (column view ids)
$ws.Range("F31").Formula = '="if(row=="&LEFT(A6,1)&" && column=="&RIGHT(A6,1)&") "&F1&".getBackground().setColorFilter(Color.parseColor("&CHAR(34)&"#ECEFF1"&CHAR(34)&"), PorterDuff.Mode.ADD);"'
$ws.Range("G31").Formula = '="if(row=="&LEFT(B6,1)&" && column=="&RIGHT(B6,1)&") "&G1&".getBackground().setColorFilter(Color.parseColor("&CHAR(34)&"#ECEFF1"&CHAR(34)&"), PorterDuff.Mode.ADD);"'
$ws.Range("H31").Formula = '="if(row=="&LEFT(C6,1)&" && column=="&RIGHT(C6,1)&") "&H1&".getBackground().setColorFilter(Color.parseColor("&CHAR(34)&"#ECEFF1"&CHAR(34)&"), PorterDuff.Mode.ADD);"'
$ws.Range("I31").Formula = '="if(row=="&LEFT(D6,1)&" && column=="&RIGHT(D6,1)&") "&I1&".getBackground().setColorFilter(Color.parseColor("&CHAR(34)&"#ECEFF1"&CHAR(34)&"), PorterDuff.Mode.ADD);"'

# Row 32 references row 7 / row 2
$ws.Range("F32").Formula = '="if(row=="&LEFT(A7,1)&" && column=="&RIGHT(A7,1)&") "&F2&".getBackground().setColorFilter(Color.parseColor("&CHAR(34)&"#ECEFF1"&CHAR(34)&"), PorterDuff.Mode.ADD);"'
$ws.Range("G32").Formula = '="if(row=="&LEFT(B7,1)&" && column=="&RIGHT(B7,1)&") "&G2&".getBackground().setColorFilter(Color.parseColor("&CHAR(34)&"#ECEFF1"&CHAR(34)&"), PorterDuff.Mode.ADD);"'
$ws.Range("H32").Formula = '="if(row=="&LEFT(C7,1)&" && column=="&RIGHT(C7,1)&") "&H2&".getBackground().setColorFilter(Color.parseColor("&CHAR(34)&"#ECEFF1"&CHAR(34)&"), PorterDuff.Mode.ADD);"'
$ws.Range("I32").Formula = '="if(row=="&LEFT(D7,1)&" && column=="&RIGHT(D7,1)&") "&I2&".getBackground().setColorFilter(Color.parseColor("&CHAR(34)&"#ECEFF1"&CHAR(34)&"), PorterDuff.Mode.ADD);"'

# Row 33 references row 8 / row 3
$ws.Range("F33").Formula = '="if(row=="&LEFT(A8,1)&" && column=="&RIGHT(A8,1)&") "&F3&".getBackground().setColorFilter(Color.parseColor("&CHAR(34)&"#ECEFF1"&CHAR(34)&"), PorterDuff.Mode.ADD);"'
$ws.Range("G33").Formula = '="if(row=="&LEFT(B8,1)&" && column=="&RIGHT(B8,1)&") "&G3&".getBackground().setColorFilter(Color.parseColor("&CHAR(34)&"#ECEFF1"&CHAR(34)&"), PorterDuff.Mode.ADD);"'
$ws.Range("H33").Formula = '="if(row=="&LEFT(C8,1)&" && column=="&RIGHT(C8,1)&") "&H3&".getBackground().setColorFilter(Color.parseColor("&CHAR(34)&"#ECEFF1"&CHAR(34)&"), PorterDuff.Mode.ADD);"'
$ws.Range("I33").Formula = '="if(row=="&LEFT(D8,1)&" && column=="&RIGHT(D8,1)&") "&I3&".getBackground().setColorFilter(Color.parseColor("&CHAR(34)&"#ECEFF1"&CHAR(34)&"), PorterDuff.Mode.ADD);"'

# Row 34 references row 9 / row 4
$ws.Range("F34").Formula = '="if(row=="&LEFT(A9,1)&" && column=="&RIGHT(A9,1)&") "&F4&".getBackground().setColorFilter(Color.parseColor("&CHAR(34)&"#ECEFF1"&CHAR(34)&"), PorterDuff.Mode.ADD);"'
$ws.Range("G34").Formula = '="if(row=="&LEFT(B9,1)&" && column=="&RIGHT(B9,1)&") "&G4&".getBackground().setColorFilter(Color.parseColor("&CHAR(34)&"#ECEFF1"&CHAR(34)&"), PorterDuff.Mode.ADD);"'
$ws.Range("H34").Formula = '="if(row=="&LEFT(C9,1)&" && column=="&RIGHT(C9,1)&") "&H4&".getBackground().setColorFilter(Color.parseColor("&CHAR(34)&"#ECEFF1"&CHAR(34)&"), PorterDuff.Mode.ADD);"'
$ws.Range("I34").Formula = '="if(row=="&LEFT(D9,1)&" && column=="&RIGHT(D9,1)&") "&I4&".getBackground().setColorFilter(Color.parseColor("&CHAR(34)&"#ECEFF1"&CHAR(34)&"), PorterDuff.Mode.ADD);"'

# --- Selection moves to match the new working range ---
$ws.Range("F31:I34").Select()
$excel.ActiveWindow.RangeSelection.Item(1).Activate()
$ws.Range("I31").Activate()
